$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117; this shifts rows 117:173 down to 118:174
$ws.Rows.Item(117).Insert()

# Populate the new row 117 with the new data record
$ws.Cells.Item(117, 1).Value = 8
$ws.Cells.Item(117, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = 45016
$ws.Cells.Item(117, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(117, 5).Value = 4
$ws.Cells.Item(117, 6).Value = 100112052
$ws.Cells.Item(117, 7).Value = "Albahaca"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 1060
$ws.Cells.Item(117, 11).Value = 2500
$ws.Cells.Item(117, 12).Value = 3000
$ws.Cells.Item(117, 13).Value = 2750
$ws.Cells.Item(117, 14).Value = "`$/docena de matas"
$ws.Cells.Item(117, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(117, 16).Value = 458
$ws.Cells.Item(117, 17).Value = 6
$ws.Cells.Item(117, 18).Value = "Hortaliza"
